$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 225
$ws1.Range("F5").Value = 2724
$ws1.Range("F7").Value = 371
$ws1.Range("F9").Value = 965

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 225
$ws4.Range("F5").Value = 2724
$ws4.Range("F7").Value = 371
$ws4.Range("F10").Value = 965
